# Trade #3 closed at 2026-02-17 13:07:34 - unknown UNKNOWN +0.000%
#
# Updates the live trading results workbook to reflect that Trade #3
# (MarketMaking strategy) has now been closed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("B3").Value = 1199.56    # Current Capital
$summary.Range("B4").Value = -0.44      # Total P&L $
$summary.Range("B5").Value = -2.93      # Total P&L %
$summary.Range("B6").Value = 3          # Total Trades
$summary.Range("B8").Value = 3          # Losing Trades

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")

$status.Range("C4").Value = 99.56       # Capital
$status.Range("D4").Value = 3           # Trades
$status.Range("E4").Value = -0.44       # P&L $
$status.Range("F4").Value = -0.44       # P&L %

# ---------------------------------------------------------------------
# All Trades sheet and MarketMaking sheet both track the same trade
# log for Trade #3 (row 4) - it is now closed.
# ---------------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("G4").Value = 0.24            # Exit Price
    $ws.Range("H4").Value = "CLOSED"        # Status
    $ws.Range("I4").Value = -55.5556        # P&L %
    $ws.Range("J4").Value = -0.3            # P&L $
    $ws.Range("K4").Value = 99.56           # Capital After
    $ws.Range("P4").Value = "early_exit"    # Exit Reason
    $ws.Range("Q4").Value = 5.51            # Duration (min)
}
